$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column C (Background Processing) to FALSE for every even data row (2,4,...,38)
for ($r = 2; $r -le 38; $r += 2) {
    $ws.Cells.Item($r, 3).Value = $false
}

# Update the selected range to match the edited column
$ws.Range("C2:C38").Select()
